# The "prediction" score column (B2:B8) was recomputed and the sheet
# re-exported; update the per-bin prediction values to the new results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2673.3021260800278
$ws.Range("B3").Value = 2003.706563849601
$ws.Range("B4").Value = 2027.0921268103539
$ws.Range("B5").Value = 2147.9399552252121
$ws.Range("B6").Value = 1920.1576005446896
$ws.Range("B7").Value = 1783.4201048542939
$ws.Range("B8").Value = 1978.9273287045266
